# Sprint 2 burndown update
# Applies the "Day 3 / Day 4" effort re-allocations that shift the
# burndown numbers (D21/F21/G21/H21) and re-points the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Display a combination of search results in one list" ---
# Remove the 2 hours that had been logged against Day 4 (G3).
$ws.Range("G3").ClearContents()

# --- Row 6: "Develop website interface..." ---
# The highlighted (colour-filled) Day-1/Day-3 cells lose their fill,
# becoming plain cells again (value in B6 / E6 is unchanged).
$ws.Range("B6").ClearFormats()
$ws.Range("E6").ClearFormats()

# --- Row 7: "Merge HTML front end... order by price, distance" ---
# Same fill clean-up on B7 / E7 (both keep their existing values/blank).
$ws.Range("B7").ClearFormats()
$ws.Range("E7").ClearFormats()

# --- Row 8: "Corrolate search results with map pins" ---
# 3 hours of work logged against Day 4.
$ws.Range("G8").Value = 3

# --- Row 9: "Display additional details for search results" ---
# Effort now spread/corrected across Day 1 - Day 4. E9 also loses its
# highlight fill (same clean-up as the other colour-filled cells above).
$ws.Range("D9").Value = 1
$ws.Range("E9").ClearFormats()
$ws.Range("E9").Value = -1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3

# --- Row 14: "Test recieving data from Java function calls..." ---
$ws.Range("D14").ClearFormats()

# --- Row 16: "Create an option to use the device's location..." ---
$ws.Range("G16").Value = 1

# --- Row 19 / Row 20: formatting-only highlight cells revert to plain ---
$ws.Range("B19").ClearFormats()
$ws.Range("F20").ClearFormats()

# --- Selection moves from G9 to B19 to match where editing finished ---
$ws.Range("B19").Select()
